$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = -0.08368629757078778
$ws.Range("C2").Value = 0.3438442104993772
$ws.Range("D2").Value = 0.4242218510161214
$ws.Range("E2").Value = 0.6513231540611168
$ws.Range("F2").Value = 0.6522263476946283
$ws.Range("B3").Value = 0.205312877006492
$ws.Range("C3").Value = 0.4985604499272415
$ws.Range("D3").Value = 0.4677225124169028
$ws.Range("E3").Value = 0.6839024143961643
$ws.Range("F3").Value = 0.6588478714022412
$ws.Range("B4").Value = 0.0422407319769525
$ws.Range("C4").Value = 0.6156476625237515
$ws.Range("D4").Value = 0.8090695585645494
$ws.Range("E4").Value = 0.8994829395628077
$ws.Range("F4").Value = 0.9076125213463587
$ws.Range("B5").Value = 0.1726989577612019
$ws.Range("C5").Value = 0.4859904821377091
$ws.Range("D5").Value = 0.5789570236283166
$ws.Range("E5").Value = 0.7608922549404197
$ws.Range("F5").Value = 0.7487137941609299
$ws.Range("B6").Value = -0.05426476430528045
$ws.Range("C6").Value = 0.5519219891770294
$ws.Range("D6").Value = 0.6681402045837443
$ws.Range("E6").Value = 0.8173984368615739
$ws.Range("F6").Value = 0.8242260833855453
$ws.Range("B7").Value = 0.0383632470686223
$ws.Range("C7").Value = 0.443900067936677
$ws.Range("D7").Value = 0.490699680860556
$ws.Range("E7").Value = 0.7004995937618779
$ws.Range("F7").Value = 0.7085918159736068
$ws.Range("B8").Value = -0.0001125680424399982
$ws.Range("C8").Value = 0.5144441553180874
$ws.Range("D8").Value = 0.5897234012367705
$ws.Range("E8").Value = 0.7679345032206656
$ws.Range("F8").Value = 0.7782428017825981
$ws.Range("B9").Value = 0.02398440652176793
$ws.Range("C9").Value = 0.5322186494472522
$ws.Range("D9").Value = 0.7385541288137565
$ws.Range("E9").Value = 0.8593917202380743
$ws.Range("F9").Value = 0.880271447287956
$ws.Range("B10").Value = -0.2556182063508592
$ws.Range("C10").Value = 0.3897054822362743
$ws.Range("D10").Value = 0.4087892466617125
$ws.Range("E10").Value = 0.6393662852088093
$ws.Range("F10").Value = 0.6081674939718213
